$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Summary")
$ws.Range("B1").Value = "BEALE_FINAL"
